# Automatische sync: 2025-06-17 15:57:44
$wb = $excel.ActiveWorkbook

# --- Update the "Logs" sheet: append new row 18 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(18, 1).Value = "Vragen over samenwerking"
$logs.Cells.Item(18, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(18, 3).Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Cells.Item(18, 4).Value = "Overig"
$logs.Cells.Item(18, 6).Value = "2025-06-17 14:59:58"
$logs.Cells.Item(18, 7).Value = "Nee"

# --- Extend conditional formatting ranges to include the new row ---
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))

# --- Update the "Dashboard" sheet: increment count for "Overig" category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 5
